# Regenerate save_data to use K instead of Strike#: update the G column
# (K = strikeouts) values for the game log rows on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new K (G column) value
$updates = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 3
    6  = 1
    7  = 1
    8  = 2
    10 = 0
    11 = 2
    12 = 1
    13 = 2
    14 = 4
    15 = 3
    17 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
